$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" (D5) and "Correspond Handback DateTime" (G5)
# on the zh-cn sheet with newly-generated report timestamps.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-02-16 15:35:18"
$wsZh.Range("G5").Value = "2016-02-16 15:36:14"

# Same two cells on the de-de sheet.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-02-16 15:35:37"
$wsDe.Range("G5").Value = "2016-02-16 15:36:44"
